# Scheduled runner update: refresh cached market-board pricing/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 362.0909
$ws.Range("I9").Value = 198.11111
$ws.Range("K9").Value = 198.11111
$ws.Range("M9").Value = -29.11111

$ws.Range("H40").Value = 3902.739
$ws.Range("J40").Value = 4779.909
$ws.Range("L40").Value = 4779.909
$ws.Range("N40").Value = -5129.909

$ws.Range("H98").Value = 52635056
$ws.Range("I98").Value = 55559148
$ws.Range("K98").Value = 55559148
$ws.Range("M98").Value = -55557650

$ws.Range("H99").Value = 452.57144
$ws.Range("I99").Value = 452.57144
$ws.Range("K99").Value = 1357.71432
$ws.Range("M99").Value = 140.28568

$ws.Range("H122").Value = 52635056
$ws.Range("I122").Value = 55559148
$ws.Range("K122").Value = 166677444
$ws.Range("M122").Value = -166674994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1321161.5
$ws.Range("I32").Value = 1408917.4
$ws.Range("K32").Value = 1408917.4
$ws.Range("M32").Value = -1408630.4

$ws.Range("H63").Value = 1317.4286
$ws.Range("I63").Value = 1080.25
$ws.Range("J63").Value = 1633.6666
$ws.Range("K63").Value = 1080.25
$ws.Range("L63").Value = 1633.6666
$ws.Range("M63").Value = -394.25
$ws.Range("N63").Value = -3005.6666

$ws.Range("H66").Value = 1317.4286
$ws.Range("I66").Value = 1080.25
$ws.Range("J66").Value = 1633.6666
$ws.Range("K66").Value = 5401.25
$ws.Range("L66").Value = 8168.333000000001
$ws.Range("M66").Value = -1969.25
$ws.Range("N66").Value = -15032.333

$ws.Range("H74").Value = 34275.613
$ws.Range("I74").Value = 42851.25
$ws.Range("K74").Value = 42851.25
$ws.Range("M74").Value = -41977.25

$ws.Range("H77").Value = 34275.613
$ws.Range("I77").Value = 42851.25
$ws.Range("K77").Value = 214256.25
$ws.Range("M77").Value = -209888.25

$ws.Range("H88").Value = 1538.1818
$ws.Range("I88").Value = 847.5
$ws.Range("J88").Value = 1932.8572
$ws.Range("K88").Value = 847.5
$ws.Range("L88").Value = 1932.8572
$ws.Range("M88").Value = -441.5
$ws.Range("N88").Value = -2744.8572

$ws.Range("H91").Value = 1538.1818
$ws.Range("I91").Value = 847.5
$ws.Range("J91").Value = 1932.8572
$ws.Range("K91").Value = 847.5
$ws.Range("L91").Value = 1932.8572
$ws.Range("M91").Value = 556.5
$ws.Range("N91").Value = -4740.8572

$ws.Range("H102").Value = 666.7143
$ws.Range("I102").Value = 666.7143
$ws.Range("K102").Value = 666.7143
$ws.Range("M102").Value = 955.2857

$ws.Range("H132").Value = 5431.8184
$ws.Range("I132").Value = 4399.7666
$ws.Range("K132").Value = 13199.2998
$ws.Range("M132").Value = -10669.2998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7938455.5
$ws.Range("I20").Value = 12822570
$ws.Range("J20").Value = 1768.25
$ws.Range("K20").Value = 12822570
$ws.Range("L20").Value = 1768.25
$ws.Range("M20").Value = -12822323
$ws.Range("N20").Value = -2262.25

$ws.Range("H22").Value = 256.16666
$ws.Range("I22").Value = 256.16666
$ws.Range("K22").Value = 256.16666
$ws.Range("M22").Value = -83.16665999999998

$ws.Range("H86").Value = 200205600
$ws.Range("I86").Value = 1000001
$ws.Range("J86").Value = 250006990
$ws.Range("K86").Value = 1000001
$ws.Range("L86").Value = 250006990
$ws.Range("M86").Value = -998878
$ws.Range("N86").Value = -250009236

$ws.Range("H89").Value = 200205600
$ws.Range("I89").Value = 1000001
$ws.Range("J89").Value = 250006990
$ws.Range("K89").Value = 5000005
$ws.Range("L89").Value = 1250034950
$ws.Range("M89").Value = -4994389
$ws.Range("N89").Value = -1250046182

$ws.Range("H99").Value = 10103349
$ws.Range("J99").Value = 18184348
$ws.Range("L99").Value = 18184348
$ws.Range("N99").Value = -18187344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5251.9756
$ws.Range("I31").Value = 2476.3147
$ws.Range("K31").Value = 2476.3147
$ws.Range("M31").Value = -2181.3147

$ws.Range("H34").Value = 5251.9756
$ws.Range("I34").Value = 2476.3147
$ws.Range("K34").Value = 2476.3147
$ws.Range("M34").Value = -2274.3147

$ws.Range("H52").Value = 70780
$ws.Range("J52").Value = 70780
$ws.Range("L52").Value = 70780
$ws.Range("N52").Value = -71368

$ws.Range("H105").Value = 7940502.5
$ws.Range("I105").Value = 11905920
$ws.Range("K105").Value = 11905920
$ws.Range("M105").Value = -11904173

$ws.Range("H134").Value = 4332.806
$ws.Range("I134").Value = 1854.1489
$ws.Range("J134").Value = 10157.65
$ws.Range("K134").Value = 5562.4467
$ws.Range("L134").Value = 30472.95
$ws.Range("M134").Value = -3027.4467
$ws.Range("N134").Value = -35542.95

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1698.909
$ws.Range("I80").Value = 1698.909
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1698.909
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -700.9090000000001
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 1698.909
$ws.Range("I83").Value = 1698.909
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 8494.545
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -3502.545
$ws.Range("N83").ClearContents()

$ws.Range("H100").Value = 25320.25
$ws.Range("J100").Value = 25320.25
$ws.Range("L100").Value = 25320.25
$ws.Range("N100").Value = -27484.25

$ws.Range("H122").Value = 2013747.4
$ws.Range("I122").Value = 3150537.2
$ws.Range("J122").Value = 2503.6155
$ws.Range("K122").Value = 9451611.600000001
$ws.Range("L122").Value = 7510.8465
$ws.Range("M122").Value = -9449161.600000001
$ws.Range("N122").Value = -12410.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6883.6665
$ws.Range("I7").Value = 5825.5
$ws.Range("J7").Value = 9000
$ws.Range("K7").Value = 5825.5
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = -5713.5
$ws.Range("N7").Value = -9224

$ws.Range("H75").Value = 41500
$ws.Range("J75").Value = 41500
$ws.Range("L75").Value = 41500
$ws.Range("N75").Value = -43372

$ws.Range("H78").Value = 41500
$ws.Range("J78").Value = 41500
$ws.Range("L78").Value = 124500
$ws.Range("N78").Value = -133860

$ws.Range("H82").Value = 4318.25
$ws.Range("I82").Value = 4867
$ws.Range("J82").Value = 3989
$ws.Range("K82").Value = 4867
$ws.Range("L82").Value = 3989
$ws.Range("M82").Value = -4506
$ws.Range("N82").Value = -4711

$ws.Range("H85").Value = 4318.25
$ws.Range("I85").Value = 4867
$ws.Range("J85").Value = 3989
$ws.Range("K85").Value = 4867
$ws.Range("L85").Value = 3989
$ws.Range("M85").Value = -3619
$ws.Range("N85").Value = -6485

$ws.Range("H93").Value = 4429.25
$ws.Range("I93").Value = 2298.9285
$ws.Range("J93").Value = 9400
$ws.Range("K93").Value = 2298.9285
$ws.Range("L93").Value = 9400
$ws.Range("M93").Value = -1050.9285
$ws.Range("N93").Value = -11896

$ws.Range("H100").Value = 4038.1667
$ws.Range("I100").Value = 2736
$ws.Range("K100").Value = 2736
$ws.Range("M100").Value = -2195

$ws.Range("H126").Value = 6883.6665
$ws.Range("I126").Value = 5825.5
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 17476.5
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -15006.5
$ws.Range("N126").Value = -31940

$ws.Range("H136").Value = 9612.521000000001
$ws.Range("I136").Value = 1735.2727
$ws.Range("J136").Value = 16833.334
$ws.Range("K136").Value = 5205.8181
$ws.Range("L136").Value = 50500.00199999999
$ws.Range("M136").Value = -2655.8181
$ws.Range("N136").Value = -55600.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 66734000
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 100100000
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 200200000
$ws.Range("M81").Value = -2939
$ws.Range("N81").Value = -200202122

$ws.Range("H84").Value = 66734000
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 100100000
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 1001000000
$ws.Range("M84").Value = -14696
$ws.Range("N84").Value = -1001010608

$ws.Range("H113").Value = 12111.363
$ws.Range("I113").Value = 13950.105
$ws.Range("J113").Value = 466
$ws.Range("K113").Value = 41850.315
$ws.Range("L113").Value = 1398
$ws.Range("M113").Value = -39680.315
$ws.Range("N113").Value = -5738

$ws.Range("H122").Value = 450995.56
$ws.Range("I122").Value = 4000004
$ws.Range("K122").Value = 12000012
$ws.Range("M122").Value = -11997562

$ws.Range("H132").Value = 17872408
$ws.Range("I132").Value = 21745454
$ws.Range("K132").Value = 65236362
$ws.Range("M132").Value = -65233832
